$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3438.9443
$ws.Range("I62").Value = 3033.4167
$ws.Range("J62").Value = 4250
$ws.Range("K62").Value = 3033.4167
$ws.Range("L62").Value = 4250
$ws.Range("M62").Value = -2409.4167
$ws.Range("N62").Value = -5498

$ws.Range("H65").Value = 3438.9443
$ws.Range("I65").Value = 3033.4167
$ws.Range("J65").Value = 4250
$ws.Range("K65").Value = 15167.0835
$ws.Range("L65").Value = 21250
$ws.Range("M65").Value = -12047.0835
$ws.Range("N65").Value = -27490

$ws.Range("H111").Value = 100003300
$ws.Range("I111").Value = 200006180
$ws.Range("J111").Value = 426.6
$ws.Range("K111").Value = 600018540
$ws.Range("L111").Value = 1279.8
$ws.Range("M111").Value = -600015473
$ws.Range("N111").Value = -7413.8

$ws.Range("H116").Value = 2265844.8
$ws.Range("I116").Value = 9617641
$ws.Range("J116").Value = 3753.6924
$ws.Range("K116").Value = 9617641
$ws.Range("L116").Value = 3753.6924
$ws.Range("M116").Value = -9614199
$ws.Range("N116").Value = -10637.6924

$ws.Range("H137").Value = 29657.457
$ws.Range("I137").Value = 596.5
$ws.Range("J137").Value = 113611.336
$ws.Range("K137").Value = 1789.5
$ws.Range("L137").Value = 340834.008
$ws.Range("M137").Value = 760.5
$ws.Range("N137").Value = -345934.008

$ws.Range("H138").Value = 3676.9663
$ws.Range("I138").Value = 1490.2593
$ws.Range("J138").Value = 4629.2417
$ws.Range("K138").Value = 4470.7779
$ws.Range("L138").Value = 13887.7251
$ws.Range("M138").Value = 669.2221
$ws.Range("N138").Value = -24167.7251

$ws.Range("H141").Value = 2398.8
$ws.Range("I141").Value = 2499.75
$ws.Range("J141").Value = 1995
$ws.Range("K141").Value = 7499.25
$ws.Range("L141").Value = 5985
$ws.Range("M141").Value = -2319.25
$ws.Range("N141").Value = -16345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 20000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1000
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = -741
$ws.Range("N23").Value = 0

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = 0

$ws.Range("H61").Value = 4693.933
$ws.Range("I61").Value = 1888.3334
$ws.Range("K61").Value = 1888.3334
$ws.Range("M61").Value = -1676.3334

$ws.Range("H110").Value = 888.7059
$ws.Range("I110").Value = 637.7778
$ws.Range("J110").Value = 1856.5714
$ws.Range("K110").Value = 637.7778
$ws.Range("L110").Value = 1856.5714
$ws.Range("M110").Value = 1407.2222
$ws.Range("N110").Value = -5946.5714

$ws.Range("H136").Value = 4693.933
$ws.Range("I136").Value = 1888.3334
$ws.Range("K136").Value = 5665.0002
$ws.Range("M136").Value = -3115.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2029.4584
$ws.Range("I107").Value = 1769.4375
$ws.Range("J107").Value = 2549.5
$ws.Range("K107").Value = 1769.4375
$ws.Range("L107").Value = 2549.5
$ws.Range("M107").Value = 150.5625
$ws.Range("N107").Value = -6389.5

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws.Range("H138").Value = 52137.582
$ws.Range("J138").Value = 52137.582
$ws.Range("L138").Value = 52137.582
$ws.Range("N138").Value = -62417.582

$ws.Range("H140").Value = 48386
$ws.Range("J140").Value = 48386
$ws.Range("L140").Value = 48386
$ws.Range("N140").Value = -58746

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2341.8823
$ws.Range("I58").Value = 2235.3333
$ws.Range("J58").Value = 2461.75
$ws.Range("K58").Value = 2235.3333
$ws.Range("L58").Value = 2461.75
$ws.Range("M58").Value = -2032.3333
$ws.Range("N58").Value = -2867.75

$ws.Range("H99").Value = 1649.3462
$ws.Range("I99").Value = 1372.7894
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 1372.7894
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = 125.2106000000001
$ws.Range("N99").Value = -5396

$ws.Range("H126").Value = 1649.3462
$ws.Range("I126").Value = 1372.7894
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 4118.3682
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -1648.3682
$ws.Range("N126").Value = -12140

$ws.Range("H132").Value = 2680.3225
$ws.Range("I132").Value = 1792.8422
$ws.Range("J132").Value = 4085.5
$ws.Range("K132").Value = 5378.5266
$ws.Range("L132").Value = 12256.5
$ws.Range("M132").Value = -2848.5266
$ws.Range("N132").Value = -17316.5

$ws.Range("H134").Value = 1650.7727
$ws.Range("I134").Value = 1066.2
$ws.Range("K134").Value = 3198.6
$ws.Range("M134").Value = -663.6000000000004

$ws.Range("H136").Value = 2341.8823
$ws.Range("I136").Value = 2235.3333
$ws.Range("J136").Value = 2461.75
$ws.Range("K136").Value = 6705.999899999999
$ws.Range("L136").Value = 7385.25
$ws.Range("M136").Value = -4155.999899999999
$ws.Range("N136").Value = -12485.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3631.5833
$ws.Range("I131").Value = 356
$ws.Range("J131").Value = 4891.423
$ws.Range("K131").Value = 1068
$ws.Range("L131").Value = 14674.269
$ws.Range("M131").Value = 3972
$ws.Range("N131").Value = -24754.269

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2734.96
$ws.Range("I132").Value = 2359
$ws.Range("K132").Value = 7077
$ws.Range("M132").Value = -4547

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7020.7646
$ws.Range("I61").Value = 7334.5625
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 7334.5625
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -7132.5625
$ws.Range("N61").Value = -2404

$ws.Range("H113").Value = 7020.7646
$ws.Range("I113").Value = 7334.5625
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 7334.5625
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -5164.5625
$ws.Range("N113").Value = -6340

$ws.Range("H128").Value = 39957.4
$ws.Range("J128").Value = 39957.4
$ws.Range("L128").Value = 39957.4
$ws.Range("N128").Value = -49917.4

$ws.Range("H132").Value = 7908.93
$ws.Range("I132").Value = 9468.906
$ws.Range("J132").Value = 3370.818
$ws.Range("K132").Value = 28406.718
$ws.Range("L132").Value = 10112.454
$ws.Range("M132").Value = -25876.718
$ws.Range("N132").Value = -15172.454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 337.85715
$ws.Range("I113").Value = 304.75
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 914.25
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1255.75
$ws.Range("N113").Value = -7340

$ws.Range("H126").Value = 1485.4
$ws.Range("I126").Value = 1501.762
$ws.Range("J126").Value = 1447.2222
$ws.Range("K126").Value = 4505.286
$ws.Range("L126").Value = 4341.6666
$ws.Range("M126").Value = -2035.286
$ws.Range("N126").Value = -9281.6666

$ws.Range("H136").Value = 3725.913
$ws.Range("I136").Value = 4203.7812
$ws.Range("J136").Value = 2633.6428
$ws.Range("K136").Value = 12611.3436
$ws.Range("L136").Value = 7900.928400000001
$ws.Range("M136").Value = -10061.3436
$ws.Range("N136").Value = -13000.9284
